$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column I ("wij") first, then column G ("wj"), so that column
# letters referenced below are still valid at the moment of deletion.
$ws.Range("I1:I4").EntireColumn.Delete() | Out-Null
$ws.Range("G1:G4").EntireColumn.Delete() | Out-Null

# After the two column deletions the sheet now has columns A:K.
# Update the dimension/shape is handled automatically by the engine;
# just make sure the remaining values match the target values exactly
# (a couple of the downstream numbers were recomputed upstream, not merely
# shifted, so set them explicitly here).

$ws.Range("H2").Value = 0.7679735396567061
$ws.Range("I2").Value = 0.6981577633242783

$ws.Range("H3").Value = 0.7475156780180909
$ws.Range("I3").Value = 0.6795597072891736

$ws.Range("H4").Value = 0.5331246591925922
$ws.Range("I4").Value = 0.1938635124336699
